$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.430.18"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.037.52"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.09"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.47"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.034.97"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("E12").Value = "  +7.04%  "
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.81"
$ws.Range("E14").Value = "  +6.62%  "
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.365.22"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.539.47"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("E18").Value = "  +4.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.035.22"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.40"
$ws.Range("E20").Value = "  +18.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.26"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.38"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.19"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.80"
$ws.Range("E25").Value = "  +4.64%  "
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.22"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("E32").Value = "  +7.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0993"
$ws.Range("E33").Value = "  -4.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.34"
$ws.Range("E34").Value = "  +4.15%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.995"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.87"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.54"
$ws.Range("E38").Value = "  +10.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.07"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.57"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.62"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "380.01"
$ws.Range("E46").Value = "  -4.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.712.37"
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.80"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.53"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("E51").Value = "  +3.83%  "

Write-Host "Applied all changes"
